# Add a "time_taken" metadata column (F) to the panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: same bold/border/centered style as the other
# header cells (B1:E1), so copy formatting from E1 onto F1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# New data cell F2: plain text value, no special style.
$ws.Range("F2").Value = "2021-10-05 13:39:00.647965"
